$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.039.22'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.563.07'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('D5').Value = "'208.50"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = "'22.09"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '1.785.31'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = '1.559.97'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '27.029.91'
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('D17').Value = "'61.88"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').Value = "'215.88"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').Value = "'7.38"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('D22').Value = "'4.14"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.25%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').Value = "'1.94"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('D25').Value = "'153.83"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').Value = "'15.05"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('E31').Value = '  +4.32%  '
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('E33').Value = '  +3.81%  '
$ws.Range('D34').Value = '1.423.61'
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('E35').Value = '  +1.71%  '
$ws.Range('E36').Value = '  +10.14%  '
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = "'0.811"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'5.80"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = "'2.33"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'1.00"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('D45').Value = "'64.76"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').Value = '1.698.82'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('D48').Value = "'86.71"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('E51').Value = '  +0.62%  '
